$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.345.46'
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").Value = '2.222.24'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.94'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -8.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.16'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +7.57%  '
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("E9").Value = '  -3.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.64'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -8.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0908'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.23'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.60'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -8.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.01'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +10.97%  '
$ws.Range("E15").Value = '  -2.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.86'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -5.85%  '
$ws.Range("D17").Value = '2.558.11'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").Value = '2.200.85'
$ws.Range("E18").Value = '  -2.95%  '
$ws.Range("D19").Value = '42.345.29'
$ws.Range("E19").Value = '  -2.80%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000105'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.44%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.11'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.75'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  +12.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '235.53'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.92'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -8.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.38'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -7.43%  '
$ws.Range("E29").Value = '  -2.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '37.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -10.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.01'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.08'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -8.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.00'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0876'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.86%  '
$ws.Range("E35").Value = '  -2.43%  '
$ws.Range("E36").Value = '  +4.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.19'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -6.84%  '
$ws.Range("E38").Value = '  -3.65%  '
$ws.Range("E39").Value = '  -2.81%  '
$ws.Range("E40").Value = '  -4.82%  '
$ws.Range("E41").Value = '  -6.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.88'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.01%  '
$ws.Range("E43").Value = '  -5.77%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.25'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -11.39%  '
$ws.Range("E46").Value = '  -4.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.32'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -7.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.26'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.40'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.64'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '100.73'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.15%  '
